$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) Update the cached "datetimeFigureOut" date field text from 10/27/2021 to
#    10/29/2021 everywhere it is rendered: the slide master's Date Placeholder
#    and every custom (slide) layout's own Date Placeholder override.
# ---------------------------------------------------------------------------
$oldDate = "10/27/2021"
$newDate = "10/29/2021"
$ppPlaceholderDate = 16

function Update-DatePlaceholders($shapes) {
    for ($j = 1; $j -le $shapes.Count; $j++) {
        $shp = $shapes.Item($j)
        $isDatePh = $false
        try {
            if ($shp.PlaceholderFormat.Type -eq $ppPlaceholderDate) {
                $isDatePh = $true
            }
        } catch {
            $isDatePh = $false
        }
        if ($isDatePh -and $shp.HasTextFrame) {
            $tr = $shp.TextFrame.TextRange
            if ($tr.Text -eq $oldDate) {
                $tr.Text = $newDate
            }
        }
    }
}

Update-DatePlaceholders $p.SlideMaster.Shapes

for ($i = 1; $i -le $p.SlideMaster.CustomLayouts.Count; $i++) {
    $layout = $p.SlideMaster.CustomLayouts.Item($i)
    Update-DatePlaceholders $layout.Shapes
}

# ---------------------------------------------------------------------------
# 2) Nudge the rotated "Freeform 34" connector shape (id 35) on slide 6 to its
#    new position/size. The shape is rotated 90 degrees (rot=5400000); Shape
#    Left/Top/Width/Height map directly onto the unrotated a:off / a:ext
#    bounding box, same as PowerPoint's real object model.
# ---------------------------------------------------------------------------
$slide = $p.Slides.Item(6)
for ($k = 1; $k -le $slide.Shapes.Count; $k++) {
    $shp = $slide.Shapes.Item($k)
    if ($shp.Id -eq 35 -and $shp.Name -eq "Freeform 34") {
        # Target EMU: off x=5579515 y=3796281, ext cx=770251 cy=1329515
        # (values below are chosen so the engine's internal point<->EMU
        # round-trip lands exactly on the target EMU integers)
        $shp.Left = 439.33192443847656
        $shp.Top = 298.9197998046875
        $shp.Width = 60.64972496032715
        $shp.Height = 104.68626022338867
        break
    }
}
